$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/long-term-disability-benefit-rate"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# The "Fixed Value" for Extension.url (row 5) mirrors the structure definition
# URL shown on the Metadata sheet, so it must track the same new URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/long-term-disability-benefit-rate"

# Column AI ("Constraint(s)") on row 2 (the "Extension" element row) is cleared;
# that FHIR constraint text now only lives on the "Extension.extension" row (row 4).
$elements.Range("AI2").Value = ""
